# Applies the "Update test data / Fix train and prepare data / Add save Excel" edit:
#  - prepare: append 3 new rows (Android Pay/PRODUCT, Spotify/ORG, Google Maps/PRODUCT)
#  - train:   drop the "span" column (column D), shifting "entity" from E to D
#  - config:  append a new "train_autosave" setting row (FALSE)
#  - make "prepare" the active sheet/tab (was "train")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. prepare: add rows for the remaining phrase-matcher entities
# ---------------------------------------------------------------------------
$prepare = $wb.Worksheets.Item("prepare")

$prepare.Range("A4").Value2 = 3
$prepare.Range("B4").Value2 = "phrase"
$prepare.Range("C4").Value2 = "Android Pay"
$prepare.Range("D4").Value2 = "PRODUCT"

$prepare.Range("A5").Value2 = 4
$prepare.Range("B5").Value2 = "phrase"
$prepare.Range("C5").Value2 = "Spotify"
$prepare.Range("D5").Value2 = "ORG"

$prepare.Range("A6").Value2 = 5
$prepare.Range("B6").Value2 = "phrase"
$prepare.Range("C6").Value2 = "Google Maps"
$prepare.Range("D6").Value2 = "PRODUCT"

# ---------------------------------------------------------------------------
# 2. train: remove the "span" column (D); "entity" shifts from E to D
# ---------------------------------------------------------------------------
$train = $wb.Worksheets.Item("train")
$train.Columns("D").Delete()

# ---------------------------------------------------------------------------
# 3. config: record the new train_autosave flag (disabled)
# ---------------------------------------------------------------------------
$config = $wb.Worksheets.Item("config")
$config.Range("A6").Value2 = "train_autosave"
$config.Range("B6").Value2 = $false

# ---------------------------------------------------------------------------
# 4. view state: "prepare" becomes the active sheet; select train!C6:D6
# ---------------------------------------------------------------------------
$train.Range("C6:D6").Select()
$prepare.Activate()
